$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ D = newPrice; E = newVolume } (only keys present are updated)
$updates = @{
    2  = @{ D = '299.70'; E = '-0.43%' }
    3  = @{ D = '31.74';  E = '1.13%' }
    4  = @{ D = '5.157';  E = '1.22%' }
    5  = @{ D = '0.08123'; E = '10.16%' }
    6  = @{ D = '2.442';  E = '2.56%' }
    7  = @{ D = '7.828';  E = '-1.67%' }
    8  = @{ D = '3.876';  E = '2.33%' }
    9  = @{ D = '0.9253'; E = '0.84%' }
    10 = @{ D = '0.1763'; E = '3.01%' }
    11 = @{ D = '0.07385'; E = '-2.53%' }
    12 = @{ D = '0.08916'; E = '10.74%' }
    13 = @{ D = '0.03032'; E = '0.57%' }
    14 = @{ E = '0.81%' }
    15 = @{ D = '0.001512'; E = '1.00%' }
    16 = @{ D = '0.005920'; E = '-4.13%' }
    17 = @{ D = '3.595';  E = '4.06%' }
    18 = @{ D = '2.285';  E = '2.72%' }
    19 = @{ E = '-0.99%' }
    20 = @{ E = '0.24%' }
    21 = @{ D = '4.097';  E = '-11.85%' }
    23 = @{ D = '0.04635'; E = '-0.20%' }
    24 = @{ D = '0.001251'; E = '2.02%' }
    25 = @{ D = '0.004547'; E = '1.27%' }
    26 = @{ D = '0.0001203'; E = '-7.32%' }
    27 = @{ D = '0.0003413'; E = '27.95%' }
    39 = @{ D = '0.01750'; E = '0.63%' }
    40 = @{ D = '0.04603'; E = '1.74%' }
    41 = @{ D = '0.006893'; E = '-4.88%' }
    42 = @{ D = '0.1379'; E = '2.31%' }
    43 = @{ D = '0.002216'; E = '-0.49%' }
    44 = @{ D = '0.009853' }
    45 = @{ D = '0.00006233'; E = '-1.06%' }
    46 = @{ E = '0.06%' }
    47 = @{ D = '0.008406'; E = '-15.87%' }
    48 = @{ D = '0.8009'; E = '-2.46%' }
    49 = @{ E = '0.06%' }
    50 = @{ E = '0.13%' }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        # Force text so Excel doesn't reinterpret numeric-looking / percent-looking
        # strings as numbers - matches the inline-string cell type in the sheet.
        $ws.Range($addr).NumberFormat = "@"
        $ws.Range($addr).Value = $cols[$col]
    }
}
